$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph right after the title ("Play Bounty Belles Free
#    & Read Our Game Review") containing the bold label "Meta description"
#    followed by the (non-bold) meta-description text.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

# The paragraph we just created is now paragraph #2 (currently empty).
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Discover the Wild West with Bounty Belles, a thrilling slot game with free spins and three jackpots. Play now for free and read our review.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicate bold
#    "Play Bounty Belles Free & Read Our Game Review" paragraph entirely
#    (search from the bottom so the paragraph we just inserted up top is
#    never touched), and turn the italic paragraph that used to follow it
#    into the new feature-image prompt text.
# ---------------------------------------------------------------------------

$oldBoldText = "Play Bounty Belles Free & Read Our Game Review"

for ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $oldBoldText) {
        $para.Range.Delete()
        break
    }
}

$oldItalicText = "Discover the Wild West with Bounty Belles, a thrilling slot game with free spins and three jackpots. Play now for free and read our review."
$newItalicText = "Create an eye-catching feature image for the game Bounty Belles that fits with the Western theme and includes a happy Maya warrior wearing glasses. The image should be in a cartoon style and draw attention to the three determined girls who are the main characters of the game. The background should feature a desert landscape with a saloon, and the game logo and jackpot should be prominently displayed. Use colors that match the Western color scheme, such as brown, gold, and red. The image should convey the thrill of the Wild Bounty feature and the excitement of winning one of the three jackpots. Make sure to include the Maya warrior in a prominent position to give a unique twist to the Western theme."

# Scope the Find/Replace to the final paragraph only, so the identical
# sentence that now also lives inside the new meta-description paragraph
# (inserted in step 1) is left untouched.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldItalicText, $true, $false, $false, $false, $false, `
                              $true, 1, $false, $newItalicText, 2) | Out-Null
